# Generate Report for Handback
# ------------------------------------------------------------------
# This models the "handback" report-generation step of the localization
# pipeline: once a target file has come back in sync with en-US, the
# Overview/status sheets get their status text updated and the per-language
# detail sheets get the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns populated (with a hyperlink on the
# target-file cell), for both data rows.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdTargetUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f13f73e447150e38bdc78c54596188675e480cd0/e2e/b7800bc2-0cf0-45b8-a082-e05e92a065ce.md"
$mdTargetName = "b7800bc2-0cf0-45b8-a082-e05e92a065ce.md"

$zhHandbackFile = "b7800bc2-0cf0-45b8-a082-e05e92a065ce.1fe39b7d06295cdf31c97f15f73fd668da8ca1be.zh-cn.xlf"
$deHandbackFile = "b7800bc2-0cf0-45b8-a082-e05e92a065ce.1fe39b7d06295cdf31c97f15f73fd668da8ca1be.de-de.xlf"

$zhHandbackDateTime = "2016-11-29 03:17:35"
$deHandbackDateTime = "2016-11-29 03:17:53"

# ---- Overview sheet: status text for both language columns ----------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311

# ---- Per-language detail sheets --------------------------------------
function Update-LangSheet($ws, $handbackFile, $handbackDateTime) {
    # Status column (C) for both rows
    $ws.Range("C2").Value = $statusNew
    $ws.Range("C3").Value = $statusNew

    # Latest Target File / Latest Handback File / Latest Handback DateTime
    $ws.Range("I2").Value = $mdTargetName
    $ws.Range("J2").Value = $handbackFile
    $ws.Range("K2").Value = $handbackDateTime

    $ws.Range("I3").Value = $mdTargetName
    $ws.Range("J3").Value = $handbackFile
    $ws.Range("K3").Value = $handbackDateTime

    # Widen columns to fit the newly-populated long values
    $ws.Columns.Item(3).ColumnWidth = 29.9777050018311
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40

    # Rebuild the hyperlinks collection so the new Target File link for I2/I3
    # lands between the existing A2 and A3 hyperlinks (matching row order),
    # then restore the HyperLink look on the existing A-column links.
    $linkA2 = $ws.Hyperlinks.Item(1)
    $a2Address = $linkA2.Address
    $a2Display = $linkA2.TextToDisplay
    $linkA3 = $ws.Hyperlinks.Item(2)
    $a3Address = $linkA3.Address
    $a3Display = $linkA3.TextToDisplay

    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $a2Address, "", "", $a2Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdTargetUrl, "", "", $mdTargetName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $a3Address, "", "", $a3Display) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdTargetUrl, "", "", $mdTargetName) | Out-Null

    $ws.Range("A2").Style = "HyperLink"
    $ws.Range("A3").Style = "HyperLink"
    $ws.Range("I2").Style = "HyperLink"
    $ws.Range("I3").Style = "HyperLink"
}

$wsZh = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $wsZh $zhHandbackFile $zhHandbackDateTime

$wsDe = $wb.Worksheets.Item("de-de")
Update-LangSheet $wsDe $deHandbackFile $deHandbackDateTime
